# Finalized Experiments with Participant Generation
# Renames sheets and updates stimulus-file / value references to new
# participant-generation timestamps.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (order preserved: GNG, NB, RS, TOL, vSAT) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1650291156809723"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16502911604610233"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16502911604630256"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16502911605100253"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16502911606030235"

# --- Sheet 1 (GNG) column B updates ---
$ws1.Range("B2").Value = "go_stims-16502911567645452.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911567797503.csv"
$ws1.Range("B4").Value = "go_stims-16502911567817419.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911568077176.csv"

# --- Sheet 2 (NB) column B updates ---
$ws2.Range("B2").Value = "ZB-match_5-16502911579900258.csv"
$ws2.Range("B3").Value = "ZB-match_6-16502911568129282.csv"
$ws2.Range("B4").Value = "TB-16502911599551232.csv"
$ws2.Range("B5").Value = "ZB-match_7-16502911578560324.csv"
$ws2.Range("B6").Value = "OB-16502911587200234.csv"
$ws2.Range("B7").Value = "TB-16502911604410257.csv"
$ws2.Range("B8").Value = "OB-16502911588310587.csv"
$ws2.Range("B9").Value = "OB-16502911585320244.csv"
$ws2.Range("B10").Value = "TB-16502911590260613.csv"

# --- Sheet 3 (RS) column B updates ---
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4 (TOL) column B updates ---
$ws4.Range("B2").Value = "MM_stims-16502911604770555.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911604650264.csv"
$ws4.Range("B4").Value = "MM_stims-165029116049303.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911604770555.csv"
$ws4.Range("B6").Value = "MM_stims-16502911605090258.csv"
$ws4.Range("B7").Value = "ZM_stims-1650291160494027.csv"

# --- Sheet 5 (vSAT) column B updates ---
$ws5.Range("B2").Value = "SAT_stims-16502911605150259.csv"
$ws5.Range("B3").Value = "SAT_stims-16502911605410266.csv"
$ws5.Range("B4").Value = "vSAT_stims-1650291160589081.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502911605570295.csv"
